$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the old "500n / C25" BOM line (row 4). Everything below it
#    shifts up one row, which realigns the remaining rows with the rest of
#    the existing (unchanged) BOM data.
# ---------------------------------------------------------------------------
$ws.Rows(4).Delete()

# ---------------------------------------------------------------------------
# 2) Small text correction on the "100n" line's Designator list (now row 3):
#    drop the redundant "C" in the range spans.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "C24, C30-33, C38-46, C48"

# ---------------------------------------------------------------------------
# 3) Append the new CPL / BOM rows (10-21) describing the USB connectors,
#    resistors, switches and the new 3-pin SWD-related parts (TJA1051T,
#    LSM6DSOX, AMS1117-3.3, oscillator, ...), replacing the old 10-pin
#    JTAG/SWD header with the simplified 3-pin SWD.
# ---------------------------------------------------------------------------

# Row 10: USB_OTG / J8
$ws.Range("A10").Value = "USB_OTG"
$ws.Range("B10").Value = "J8"

# Row 11: PWR_USB /   J9
$ws.Range("A11").Value = "PWR_USB"
$ws.Range("B11").Value = "  J9"

# Row 12: 10k / R1 / 603 / C191124
$ws.Range("A12").Value = "10k"
$ws.Range("B12").Value = "R1"
$ws.Range("C12").Value = 603
$ws.Range("D12").Value = "C191124"

# Row 13: 120 / " R4, R5" / 603 / C177673
$ws.Range("A13").Value = 120
$ws.Range("B13").Value = " R4, R5"
$ws.Range("C13").Value = 603
$ws.Range("D13").Value = "C177673"

# Row 14: 0 / "    R3" / 603
$ws.Range("A14").Value = 0
$ws.Range("B14").Value = "    R3"
$ws.Range("C14").Value = 603

# Row 15: SW_Push / "    SW1" / SPST_PTS645
$ws.Range("A15").Value = "SW_Push"
$ws.Range("B15").Value = "    SW1"
$ws.Range("C15").Value = "SPST_PTS645"

# Row 16: "    SW2" / SW_SPDT / PCM12
$ws.Range("A16").Value = "    SW2"
$ws.Range("B16").Value = "SW_SPDT"
$ws.Range("C16").Value = "PCM12"

# Row 17: TJA1051T / " U2, U3" / SOIC-8 / C5342108
$ws.Range("A17").Value = "TJA1051T"
$ws.Range("B17").Value = " U2, U3"
$ws.Range("C17").Value = "SOIC-8"
$ws.Range("D17").Value = "C5342108"

# Row 18: LSM6DSOX / "    U4" / LGA-14 / C481766
$ws.Range("A18").Value = "LSM6DSOX"
$ws.Range("B18").Value = "    U4"
$ws.Range("C18").Value = "LGA-14"
$ws.Range("D18").Value = "C481766"

# Row 19: AMS1117-3.3 / "    U5" / SOT-223-3 / C6186
$ws.Range("A19").Value = "AMS1117-3.3"
$ws.Range("B19").Value = "    U5"
$ws.Range("C19").Value = "SOT-223-3"
$ws.Range("D19").Value = "C6186"

# Row 20: (no comment) / "    U6" / SOT-23 / C135998
$ws.Range("B20").Value = "    U6"
$ws.Range("C20").Value = "SOT-23"
$ws.Range("D20").Value = "C135998"

# Row 21: 16MHz / "    Y1" / SMD_3225
$ws.Range("A21").Value = "16MHz"
$ws.Range("B21").Value = "    Y1"
$ws.Range("C21").Value = "SMD_3225"

# ---------------------------------------------------------------------------
# 4) Match the "JLCPCB PART#" styling (the Microsoft YaHei / dark-grey font
#    already used for D3:D9) on the newly-added part-number cells, and bump
#    the row height on those rows the same way the existing ones are.
# ---------------------------------------------------------------------------
$ws.Range("D9").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows(12).RowHeight = 16.5
$ws.Rows(13).RowHeight = 16.5
$ws.Rows(17).RowHeight = 16.5
$ws.Rows(18).RowHeight = 16.5
$ws.Rows(20).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 5) Restore the selection to B2, matching the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("B2").Select()
